$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" (Changed) date in column C for rows 2-9
# from 2023-10-22 (serial 45221) to 2023-10-25 (serial 45224)
$newDate = Get-Date -Year 2023 -Month 10 -Day 25 -Hour 0 -Minute 0 -Second 0 -Millisecond 0

for ($row = 2; $row -le 9; $row++) {
    $ws.Cells.Item($row, 3).Value = $newDate
}
